# Generate Report for Handoff
#
# The localization run finished: every "In Translation" status becomes
# "Ready for handoff", and the handoff/generate timestamps advance to the
# moment the report was produced. Excel auto-widens the Status-ish columns
# to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$overview.Range("E2").Value = "Ready for handoff"   # Overview!zh-cn
$overview.Range("F2").Value = "Ready for handoff"   # Overview!de-de
$zhcn.Range("C2").Value = "Ready for handoff"        # zh-cn!Status
$dede.Range("C2").Value = "Ready for handoff"        # de-de!Status

# --- Timestamps advance to the handoff-generation run -----------------------
$overview.Range("G2").Value = "2016-08-24 19:08:27"  # Latest HO Xliff Generate Date
$dede.Range("H2").Value     = "2016-08-24 19:08:27"  # de-de!Latest Handoff Datetime
$zhcn.Range("H2").Value     = "2016-08-24 19:08:22"  # zh-cn!Latest Handoff Datetime

# --- Column widths grow to fit "Ready for handoff" --------------------------
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
